$d = $word.ActiveDocument

# Mapping of old text -> new text (date line + 25 multiplication answers)
$replacements = @(
    @("2024-07-29 Monday", "2024-07-30 Tuesday"),
    @("105×2=210", "363×6=2178"),
    @("295×4=1180", "264×2=528"),
    @("133×7=931", "691×6=4146"),
    @("803×6=4818", "940×9=8460"),
    @("383×5=1915", "742×4=2968"),
    @("644×7=4508", "485×8=3880"),
    @("599×4=2396", "845×8=6760"),
    @("298×5=1490", "660×5=3300"),
    @("400×3=1200", "813×9=7317"),
    @("646×7=4522", "811×9=7299"),
    @("635×4=2540", "628×4=2512"),
    @("543×3=1629", "574×4=2296"),
    @("496×6=2976", "218×7=1526"),
    @("425×6=2550", "941×7=6587"),
    @("704×5=3520", "893×6=5358"),
    @("831×3=2493", "869×2=1738"),
    @("157×8=1256", "615×6=3690"),
    @("151×2=302", "121×5=605"),
    @("920×4=3680", "300×2=600"),
    @("295×9=2655", "199×8=1592"),
    @("246×7=1722", "884×8=7072"),
    @("999×7=6993", "796×5=3980"),
    @("595×4=2380", "592×2=1184"),
    @("987×2=1974", "931×9=8379"),
    @("495×5=2475", "447×2=894")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
